# Applies:
#  - Slide 2 ("Outline"): split "SBS requirements" / "Fastbus Readout" into
#    two runs each, and insert two new red-italic sub-bullets ("Data Event
#    Flow" and "Event Switching") underneath them; turn on shrink-text
#    autofit for the placeholder.
#  - Slide 9: reposition/resize the screenshot picture, delete two of the
#    three "L2" annotation textboxes, and reposition/rename the remaining one.

$p = $ppt.ActivePresentation

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Slide 2 - Outline
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$content = Get-ShapeByName $s2 "Content Placeholder 2"
$tf = $content.TextFrame
$tr = $tf.TextRange

# Paragraph 1: "SBS requirements" -> "SBS " + "requirements" (two runs with
# identical formatting - re-assigning the substring's own text forces the
# engine to split the run without touching any formatting attribute).
$para1 = $tr.Paragraphs(1, 1)
$firstWord = $tr.Characters($para1.Start, 4)
$firstWord.Text = $firstWord.Text

# Insert the new "Data Event Flow" sub-bullet right after paragraph 1.
$para1 = $tr.Paragraphs(1, 1)
$para1.InsertAfter([char]13) | Out-Null
$newPara1 = $tr.Paragraphs(2, 1)
$newPara1.IndentLevel = 2
$newPara1.Font.Italic = -1
$newPara1.Font.Color.RGB = 255
$newPara1.Text = "Data Event Flow"

# Paragraph "Fastbus Readout" is now paragraph 3 -> split into "Fastbus" +
# " Readout".
$para3 = $tr.Paragraphs(3, 1)
$firstWord2 = $tr.Characters($para3.Start, 7)
$firstWord2.Text = $firstWord2.Text

# Insert the new "Event Switching" sub-bullet right after paragraph 3.
$para3 = $tr.Paragraphs(3, 1)
$para3.InsertAfter([char]13) | Out-Null
$newPara2 = $tr.Paragraphs(4, 1)
$newPara2.IndentLevel = 2
$newPara2.Font.Italic = -1
$newPara2.Font.Color.RGB = 255
$newPara2.Text = "Event Switching"

# Shrink text on overflow (adds <a:normAutofit/> to the placeholder bodyPr).
$tf.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 9 - DAQ configuration diagram
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)

$pic = Get-ShapeByName $s9 "Picture 2"
$pic.Left = 186520 / 12700.0
$pic.Top = 944540 / 12700.0
$pic.Width = 8652680 / 12700.0
$pic.Height = 5227660 / 12700.0

$tb20 = Get-ShapeByName $s9 "TextBox 20"
$tb20.Delete()

$tb25 = Get-ShapeByName $s9 "TextBox 25"
$tb25.Delete()

$tb26 = Get-ShapeByName $s9 "TextBox 26"
$tb26.Left = 4724400 / 12700.0
$tb26.Top = 2992045 / 12700.0
$tb26.Name = "TextBox 25"
